# ============================================================================
# Edit script: restructure PlayerPerformance_4421.xlsx
#
#   - Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling"
#     and replace the full scorecard URL values with just the numeric
#     MatchCode extracted from the URL
#   - Clear out empty INNING_NUMBER cells on "ODI Batting" (column B)
#   - Insert a new "Player Info" sheet before "ODI Batting" describing the
#     player the workbook is about
#   - Insert a new "ODI Batting Extra" sheet after "ODI Bowling" with
#     additional per-match batting detail
# ============================================================================

$wb = $excel.ActiveWorkbook

# Helper: decide whether a value needs a leading apostrophe so that Excel
# keeps storing it as text instead of silently converting it to a number,
# date or percentage.
function Needs-TextPrefix($val) {
    if ($val -eq $null -or $val -eq "") { return $false }
    if ($val -match '^-?\d+(\.\d+)?$') { return $true }
    if ($val -match '^\d{1,2}/\d{1,2}/\d{2,4}$') { return $true }
    if ($val -match '^\d+(\.\d+)?%$') { return $true }
    return $false
}

# Helper: assign a text value to a cell/range, forcing text storage when the
# value looks numeric/date/percent-like.
function Set-TextValue($range, $val) {
    if ($val -ne $null -and (Needs-TextPrefix $val)) {
        $range.Value = "'" + $val
    } else {
        $range.Value = $val
    }
}

# Helper: apply the bold / bordered / centered header look used throughout
# the workbook to a header range.
function Format-HeaderRange($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 1. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingRows = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $dAddr = $wsBatting.Range("D$r")
    $url = $dAddr.Value2
    if ($url -ne $null -and $url -ne "") {
        $code = $url -replace '.*MatchCode=', ''
        Set-TextValue $dAddr $code
    }

    $bAddr = $wsBatting.Range("B$r")
    $bVal = $bAddr.Value2
    if ($bVal -eq $null -or $bVal -eq "") {
        $bAddr.Value = ""
    }
}

# ---------------------------------------------------------------------------
# 2. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingRows = $wsBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $bAddr = $wsBowling.Range("B$r")
    $url = $bAddr.Value2
    if ($url -ne $null -and $url -ne "") {
        $code = $url -replace '.*MatchCode=', ''
        Set-TextValue $bAddr $code
    }
}

# ---------------------------------------------------------------------------
# 3. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $piHeaders.Length; $c++) {
    $wsPlayerInfo.Cells.Item(1, $c + 1).Value = $piHeaders[$c]
}
Format-HeaderRange $wsPlayerInfo.Range("A1:D1")

$piRow2 = @("4421", "Mark A Wood", "Right Handed", "Right Arm Fast")
for ($c = 0; $c -lt $piRow2.Length; $c++) {
    Set-TextValue ($wsPlayerInfo.Cells.Item(2, $c + 1)) $piRow2[$c]
}

# ---------------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet, inserted after "ODI Bowling"
# ---------------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add($null, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $exHeaders.Length; $c++) {
    $wsExtra.Cells.Item(1, $c + 1).Value = $exHeaders[$c]
}
Format-HeaderRange $wsExtra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION (numeric or blank), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exRows = @(
    @("4260", 11,   "0", "0", "",      "NO"),
    @("4297", $null, "", "", "",      "NO"),
    @("4308", 11,   "", "", "",      "NO"),
    @("4314", $null, "", "", "",      "NO"),
    @("4321", 11,   "", "", "",      "NO"),
    @("4326", $null, "", "", "",      "NO"),
    @("4331", $null, "", "", "",      "NO"),
    @("4336", 11,   "", "", "",      "NO"),
    @("4342", 10,   "1", "0", "4.35%", "NO"),
    @("4346", 10,   "", "", "",      "NO"),
    @("4354", 11,   "", "", "",      "NO"),
    @("4355", $null, "", "", "",      ""),
    @("4429", $null, "", "", "",      ""),
    @("4431", $null, "", "", "",      ""),
    @("4454", $null, "", "", "",      ""),
    @("4457", $null, "", "", "",      ""),
    @("4469", $null, "", "", "",      ""),
    @("4470", $null, "", "", "",      ""),
    @("4711", $null, "", "", "",      ""),
    @("4713", $null, "", "", "",      "")
)

$r = 2
foreach ($row in $exRows) {
    Set-TextValue ($wsExtra.Cells.Item($r, 1)) $row[0]

    $posCell = $wsExtra.Cells.Item($r, 2)
    if ($row[1] -ne $null) {
        $posCell.Value = $row[1]
    }

    Set-TextValue ($wsExtra.Cells.Item($r, 3)) $row[2]
    Set-TextValue ($wsExtra.Cells.Item($r, 4)) $row[3]
    Set-TextValue ($wsExtra.Cells.Item($r, 5)) $row[4]
    Set-TextValue ($wsExtra.Cells.Item($r, 6)) $row[5]

    $r++
}

Write-Output "Done"
